$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Strike through the whole paragraph "Меню для обычного пользователя
#    только центры и помещения." (paragraph mark + run get <w:strike/>).
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("Меню для обычного пользователя только центры и помещения.") | Out-Null
$menuPara = $findRng.Paragraphs(1)
$menuPara.Range.Font.StrikeThrough = 1

# ---------------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark from the end of the
#    "Календарь для помещений даты" paragraph down to a brand new, empty
#    paragraph at the very end of the document, inserting four new
#    paragraphs in between.
# ---------------------------------------------------------------------------

# Drop the bookmark from its current spot; we'll re-add it later.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# Locate the "Календарь для помещений даты" paragraph.
$calRng = $d.Content
$calRng.Find.Execute("Календарь для помещений") | Out-Null
$calPara = $calRng.Paragraphs(1)

# Insert four fresh empty paragraphs right after it.
$anchor = $calPara.Range
$anchor.Collapse(0)
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()
$anchor.InsertParagraphAfter()

$calIndex = $calPara.Index

# Fill in the new paragraphs with their text. The second one
# (calIndex + 2) is meant to stay empty, so give it throwaway text first
# and strip it back out below -- that avoids leaving a stray empty <w:r/>
# behind (InsertParagraphAfter seeds each new paragraph with one).
$d.Paragraphs($calIndex + 1).Range.InsertAfter("Фотографии для центров")
$d.Paragraphs($calIndex + 2).Range.InsertAfter("placeholder")
$d.Paragraphs($calIndex + 3).Range.InsertAfter("Свои центры ")
$d.Paragraphs($calIndex + 4).Range.InsertAfter("Свои помещения для лендлорда")

$emptyPara = $d.Paragraphs($calIndex + 2)
$emptyRng = $emptyPara.Range
$emptyRng.MoveEnd(1, -1) | Out-Null
$emptyRng.Text = ""

# The old "Фотографии для центров" paragraph is now pushed after these new
# ones; turn it into the bookmark-only trailing paragraph.
$lastPara = $d.Paragraphs($calIndex + 5)

# Re-add the bookmark at the start of that paragraph *before* the old text
# is removed (placing it once the paragraph is truly the last one in the
# document trips an edge case, so do it while later content still exists).
$bmTarget = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmTarget)

# Now strip the old run text, leaving just the bookmark in its own paragraph.
$lastPara2 = $d.Paragraphs($calIndex + 5)
$textRng = $lastPara2.Range
$textRng.MoveEnd(1, -1) | Out-Null
$textRng.Text = ""
